$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from row 175 into the new rows 176:188
$ws.Range("B175:G175").Copy()
$ws.Range("B176:G188").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 176
$ws.Cells.Item(176,2).Value2 = 6
$ws.Cells.Item(176,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(176,4).Value2 = 71
$ws.Cells.Item(176,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(176,6).Value2 = '
6:28
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
Criação da classe UserSS - o prefixo "SS" siginifica que será um usuário que atende o contrato/interface do Spring Security'
$ws.Cells.Item(176,7).Value2 = '
'
$ws.Rows.Item(176).RowHeight = 195

# Row 177
$ws.Cells.Item(177,2).Value2 = 6
$ws.Cells.Item(177,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(177,4).Value2 = 71
$ws.Cells.Item(177,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(177,6).Value2 = '7:52
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
INteressante - caso queiramos implementar alguma lógica na classe UserSS, como por exemplo, algum tempo de expiração para a sessão do usuário, basta implementar nos metodos exigidos pela interface UserDetails - no caso o metodo "isAccountNonExpired"'
$ws.Rows.Item(177).RowHeight = 105

# Row 178
$ws.Cells.Item(178,2).Value2 = 6
$ws.Cells.Item(178,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(178,4).Value2 = 71
$ws.Cells.Item(178,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(178,6).Value2 = '
10:06
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
criação da classe UserDetailServiceImpl - que implementa o contrato do Spring Security UserDetailsService - esta interface do Spring Security permite a busca pelo nome do usuário'
$ws.Rows.Item(178).RowHeight = 105

# Row 179
$ws.Cells.Item(179,2).Value2 = 6
$ws.Cells.Item(179,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(179,4).Value2 = 71
$ws.Cells.Item(179,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(179,6).Value2 = '
12:27
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
excessão UserNameNotFoundException do spring security por estar dentro do contexto de segurança'
$ws.Rows.Item(179).RowHeight = 90

# Row 180
$ws.Cells.Item(180,2).Value2 = 6
$ws.Cells.Item(180,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(180,4).Value2 = 71
$ws.Cells.Item(180,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(180,6).Value2 = '
13:59
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
criado construtor convertendo uma lista de perfis recebida por parametro para uma Collection<? extends GrantedAuthority> que o spring security exige'
$ws.Rows.Item(180).RowHeight = 105

# Row 181
$ws.Cells.Item(181,2).Value2 = 6
$ws.Cells.Item(181,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(181,4).Value2 = 71
$ws.Cells.Item(181,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(181,6).Value2 = '14:13
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
em aula anterior foi observado a necessidade de usar o prefixo ROLE_ na criação do enum de perfis de usuario ... nesta aula mostra o por que deste prefixo e seu devido uso'
$ws.Rows.Item(181).RowHeight = 90

# Row 182
$ws.Cells.Item(182,2).Value2 = 6
$ws.Cells.Item(182,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(182,4).Value2 = 71
$ws.Cells.Item(182,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(182,6).Value2 = '20:10
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
criação de chaves no arquivo application.properties para a palabra secreta que sera embaralhada no toke e o tempo de expiração (em milissegundos) da sessão/requisição'
$ws.Rows.Item(182).RowHeight = 90

# Row 183
$ws.Cells.Item(183,2).Value2 = 6
$ws.Cells.Item(183,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(183,4).Value2 = 71
$ws.Cells.Item(183,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(183,6).Value2 = '21:25
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
criação da classe auxiliar JWT Security'
$ws.Rows.Item(183).RowHeight = 60

# Row 184
$ws.Cells.Item(184,2).Value2 = 6
$ws.Cells.Item(184,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(184,4).Value2 = 71
$ws.Cells.Item(184,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 1)'
$ws.Cells.Item(184,6).Value2 = '23:24
6. Autenticação e autorização com tokens JWT
71. Implementando autenticacao e geracao do token JWT (PARTE 1)
criação do metodo generateToken - uso da biblioteca JWT pela primeira vez - o builder contido no retorno deste metodo é quem gera o token'
$ws.Rows.Item(184).RowHeight = 90

# Row 185
$ws.Cells.Item(185,2).Value2 = 6
$ws.Cells.Item(185,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(185,4).Value2 = 72
$ws.Cells.Item(185,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 2)'
$ws.Cells.Item(185,6).Value2 = '0:11
6. Autenticação e autorização com tokens JWT
72. Implementando autenticacao e geracao do token JWT - PARTE 2
criação de um filtro de autenticação - que intercepta a requisição, executa algo antes, e depois se der certo devolve a execução para a requisição continuar normalmente'
$ws.Cells.Item(185,7).Value2 = '





'
$ws.Rows.Item(185).RowHeight = 105

# Row 186
$ws.Cells.Item(186,2).Value2 = 6
$ws.Cells.Item(186,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(186,4).Value2 = 72
$ws.Cells.Item(186,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 2)'
$ws.Cells.Item(186,6).Value2 = '1:28
6. Autenticação e autorização com tokens JWT
72. Implementando autenticacao e geracao do token JWT - PARTE 2
para que este filtro seja um filtro de autenticação, e necessário estender para um filtro do Spring Security chamado UsernamePasswordAuthenticationFilter. Quando criamos uma classe que estende UsernamePasswordAuthenticationFilter, automaticamente o Spring Security saberá que este filtro terá que interceptar a requisição de login (endpoint /login)... inclusive esse endpoint de sufixo "/login" é padrao reservado do Spring Security tbm'
$ws.Rows.Item(186).RowHeight = 150

# Row 187
$ws.Cells.Item(187,2).Value2 = 6
$ws.Cells.Item(187,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(187,4).Value2 = 72
$ws.Cells.Item(187,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 2)'
$ws.Cells.Item(187,6).Value2 = '6:44
6. Autenticação e autorização com tokens JWT
72. Implementando autenticacao e geracao do token JWT - PARTE 2
fim da implementação do metodo que autentica usuario com o framework Spring Security'
$ws.Rows.Item(187).RowHeight = 75

# Row 188
$ws.Cells.Item(188,2).Value2 = 6
$ws.Cells.Item(188,3).Value2 = 'Autenticação e autorização com tokens JWT'
$ws.Cells.Item(188,4).Value2 = 72
$ws.Cells.Item(188,5).Value2 = 'Implementando autenticacao e geracao do token JWT (PARTE 2)'
$ws.Cells.Item(188,6).Value2 = '9:55
6. Autenticação e autorização com tokens JWT
72. Implementando autenticacao e geracao do token JWT - PARTE 2
para o teste funcionar, é necessário ter inserido a atualização mencionada na aula anterior, referente ao erro 403 quando o correto é gerar um erro 401 quando o usuario insere dados invalidos de login'
$ws.Rows.Item(188).RowHeight = 105

# Resize the Excel table to include the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:G188"))

# Update the view: selection + scroll position
$ws.Range("E187").Select()
$excel.ActiveWindow.ScrollRow = 182
$excel.ActiveWindow.ScrollColumn = 1
